$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    732621518415.6881,
    244544162085.86,
    39021256108.05149,
    33240319443.19441,
    24192102645.34322,
    13626977438.05725,
    11276561340.28382,
    9044693106.554544,
    8216551625.333233,
    7983491076.597414,
    7616338019.167827,
    7579762307.742823,
    6682877131.788206,
    6125944627.159385,
    5176081409.588452,
    5036113315.379909,
    4460651036.843834,
    3657433322.423396,
    3446382275.144342,
    3425579876.474032,
    3352688912.854328,
    3006698295.104897,
    2975587062.48001,
    2833196183.419603,
    2432555563.773227
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}

$wb.Save()
